$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 21) down to the new rows
# so the new rows inherit the same cell styles (e.g. left-aligned email column,
# boolean style on is_active) that the rest of the table uses.
$ws.Range("A21:L21").Copy()
$ws.Range("A22:L30").PasteSpecial(-4122)  # xlPasteFormats

# Row 22: user Magdalena Weber
$ws.Range("A22").Value = 110021
$ws.Range("B22").Value = 7316931025
$ws.Range("C22").Value = "Magdalena Weber"
$ws.Range("D22").Value = "magdalena.weber@xyz.com"
$ws.Range("E22").Value = 932122450
$ws.Range("F22").Value = "ACT"
$ws.Range("G22").Value = "eng"
$ws.Range("H22").Value = "PWD"
$ws.Range("I22").Value = $true
$ws.Range("J22").Value = "superadmin"
$ws.Range("K22").Value = "now()"
$ws.Range("L22").Value = "now()"

# Row 23: user Adrienne Hoffman
$ws.Range("A23").Value = 110022
$ws.Range("B23").Value = 9137847236
$ws.Range("C23").Value = "Adrienne Hoffman"
$ws.Range("D23").Value = "adrienne.hoffman@xyz.com"
$ws.Range("E23").Value = 848488000
$ws.Range("F23").Value = "ACT"
$ws.Range("G23").Value = "eng"
$ws.Range("H23").Value = "PWD"
$ws.Range("I23").Value = $true
$ws.Range("J23").Value = "superadmin"
$ws.Range("K23").Value = "now()"
$ws.Range("L23").Value = "now()"

# Row 24: user Adrienne Mcgee
$ws.Range("A24").Value = 110023
$ws.Range("B24").Value = 8428758532
$ws.Range("C24").Value = "Adrienne Mcgee"
$ws.Range("D24").Value = "adrienne.mcgee@xyz.com"
$ws.Range("E24").Value = 894773246
$ws.Range("F24").Value = "ACT"
$ws.Range("G24").Value = "eng"
$ws.Range("H24").Value = "PWD"
$ws.Range("I24").Value = $true
$ws.Range("J24").Value = "superadmin"
$ws.Range("K24").Value = "now()"
$ws.Range("L24").Value = "now()"

# Row 25: user Amare Coleman
$ws.Range("A25").Value = 110024
$ws.Range("B25").Value = 9804209494
$ws.Range("C25").Value = "Amare Coleman"
$ws.Range("D25").Value = "amare.coleman@xyz.com"
$ws.Range("E25").Value = 956554588
$ws.Range("F25").Value = "ACT"
$ws.Range("G25").Value = "eng"
$ws.Range("H25").Value = "PWD"
$ws.Range("I25").Value = $true
$ws.Range("J25").Value = "superadmin"
$ws.Range("K25").Value = "now()"
$ws.Range("L25").Value = "now()"

# Row 26: user Dawson Ibarra
$ws.Range("A26").Value = 110025
$ws.Range("B26").Value = 7105248214
$ws.Range("C26").Value = "Dawson Ibarra"
$ws.Range("D26").Value = "dawson.ibarra@xyz.com"
$ws.Range("E26").Value = 765455583
$ws.Range("F26").Value = "ACT"
$ws.Range("G26").Value = "eng"
$ws.Range("H26").Value = "PWD"
$ws.Range("I26").Value = $true
$ws.Range("J26").Value = "superadmin"
$ws.Range("K26").Value = "now()"
$ws.Range("L26").Value = "now()"

# Row 27: user Elvis Mcmillan
$ws.Range("A27").Value = 110026
$ws.Range("B27").Value = 9316557128
$ws.Range("C27").Value = "Elvis Mcmillan"
$ws.Range("D27").Value = "elvis.mcmillan@xyz.com"
$ws.Range("E27").Value = 884282274
$ws.Range("F27").Value = "ACT"
$ws.Range("G27").Value = "eng"
$ws.Range("H27").Value = "PWD"
$ws.Range("I27").Value = $true
$ws.Range("J27").Value = "superadmin"
$ws.Range("K27").Value = "now()"
$ws.Range("L27").Value = "now()"

# Row 28: user Steve George
$ws.Range("A28").Value = 110027
$ws.Range("B28").Value = 8103486949
$ws.Range("C28").Value = "Steve George"
$ws.Range("D28").Value = "steve.george@xyz.com"
$ws.Range("E28").Value = 971073663
$ws.Range("F28").Value = "ACT"
$ws.Range("G28").Value = "eng"
$ws.Range("H28").Value = "PWD"
$ws.Range("I28").Value = $true
$ws.Range("J28").Value = "superadmin"
$ws.Range("K28").Value = "now()"
$ws.Range("L28").Value = "now()"

# Row 29: user Colton Elliott
$ws.Range("A29").Value = 110028
$ws.Range("B29").Value = 9601932866
$ws.Range("C29").Value = "Colton Elliott"
$ws.Range("D29").Value = "colton.elliott@xyz.com"
$ws.Range("E29").Value = 809908673
$ws.Range("F29").Value = "ACT"
$ws.Range("G29").Value = "eng"
$ws.Range("H29").Value = "PWD"
$ws.Range("I29").Value = $true
$ws.Range("J29").Value = "superadmin"
$ws.Range("K29").Value = "now()"
$ws.Range("L29").Value = "now()"

# Row 30: user Carolyn Rodriguez
$ws.Range("A30").Value = 110029
$ws.Range("B30").Value = 9317596765
$ws.Range("C30").Value = "Carolyn Rodriguez"
$ws.Range("D30").Value = "carolyn.rodriguez@xyz.com"
$ws.Range("E30").Value = 818876429
$ws.Range("F30").Value = "ACT"
$ws.Range("G30").Value = "eng"
$ws.Range("H30").Value = "PWD"
$ws.Range("I30").Value = $true
$ws.Range("J30").Value = "superadmin"
$ws.Range("K30").Value = "now()"
$ws.Range("L30").Value = "now()"

# Match the author's final selection/scroll state: the new rows are selected
$ws.Range("A22:A30").Select()

Write-Host "done"